$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Pro16xD"
$ws.Range("B4").Value = "NGC-571/T1398 OR TC-187"
$ws.Range("B4").Style = "Normal"
$ws.Range("F8").Value = 0.319
$ws.Range("J8").Value = 0.395
$ws.Range("N8").Value = 0.395
$ws.Range("O8").Value = 0.395

$ws.Range("B8").Select()
